$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 19:52"

# Refresh the country data table (re-sorted by Casos totales desc.,
# since the underlying counts moved and Excel re-sorted the range)
$data = @(
  @("Estados Unidos", 327871, 16514, 16700, 301846, 8519, 873, 9325),
  @("España", 130759, 4591, 38080, 80261, 6861, 471, 12418),
  @("Italia", 128948, 4316, 21815, 91246, 3977, 525, 15887),
  @("Alemania", 98772, 2680, 26400, 70845, 3936, 83, 1527),
  @("Francia", 89953, 0, 15438, 66955, 6838, 0, 7560),
  @("China", 81669, 30, 76964, 1376, 295, 3, 3329),
  @("Iran", 58226, 2483, 19736, 34887, 4103, 151, 3603),
  @("Reino Unido", 47806, 5903, 135, 42737, 1559, 621, 4934),
  @("Turquia", 27069, 3135, 1042, 25453, 1381, 73, 574),
  @("Suiza", 21100, 595, 6415, 14000, 391, 19, 685),
  @("Belgica", 19691, 1260, 3751, 14493, 1261, 164, 1447),
  @("Paises Bajos", 17851, 1224, 250, 15835, 1360, 115, 1766),
  @("Canada", 14426, 514, 2603, 11565, 426, 27, 258),
  @("Austria", 12051, 270, 2998, 8849, 244, 18, 204),
  @("Portugal", 11278, 754, 75, 10908, 267, 29, 295),
  @("Brasil", 10568, 208, 127, 9986, 296, 10, 455),
  @("Corea del Sur", 10237, 81, 6463, 3591, 55, 6, 183),
  @("Israel", 8018, 167, 477, 7493, 127, 4, 48),
  @("Suecia", 6830, 387, 205, 6224, 541, 28, 401),
  @("Noruega", 5687, 137, 32, 5584, 89, 9, 71),
  @("Australia", 5687, 137, 2315, 3337, 91, 5, 35),
  @("Rusia", 5389, 658, 355, 4989, 8, 2, 45),
  @("Irlanda", 4604, 0, 25, 4442, 148, 0, 137),
  @("Chequia", 4543, 71, 96, 4380, 86, 8, 67),
  @("Chile", 4471, 310, 618, 3819, 307, 7, 34),
  @("Dinamarca", 4369, 292, 1327, 2863, 142, 18, 179),
  @("Polonia", 4102, 475, 134, 3874, 50, 15, 94),
  @("Rumania", 3864, 251, 374, 3339, 141, 5, 151),
  @("Malasia", 3662, 179, 1005, 2596, 99, 4, 61),
  @("Ecuador", 3646, 181, 100, 3366, 100, 8, 180),
  @("India", 3588, 0, 229, 3260, 0, 0, 99),
  @("Filipinas", 3246, 152, 64, 3030, 1, 8, 152),
  @("Japon", 3139, 0, 514, 2548, 64, 0, 77),
  @("Pakistan", 3123, 305, 208, 2870, 18, 4, 45),
  @("Luxemburgo", 2804, 75, 500, 2268, 33, 5, 36),
  @("Arabia Saudita", 2385, 206, 488, 1863, 41, 5, 34),
  @("Peru", 2281, 535, 914, 1284, 88, 10, 83),
  @("Indonesia", 2273, 181, 164, 1911, 0, 7, 198),
  @("Tailandia", 2169, 102, 793, 1353, 23, 3, 23),
  @("Finlandia", 1927, 45, 300, 1599, 73, 3, 28),
  @("Serbia", 1908, 284, 54, 1803, 98, 7, 51),
  @("Mexico", 1890, 202, 633, 1178, 1, 19, 79),
  @("Panama", 1801, 0, 13, 1742, 75, 0, 46),
  @("Emiratos Arabes Unidos", 1798, 293, 144, 1644, 1, 0, 10),
  @("Republica Dominicana", 1745, 167, 17, 1646, 147, 5, 82),
  @("Grecia", 1735, 62, 78, 1584, 93, 5, 73),
  @("Sudafrica", 1585, 0, 95, 1481, 7, 0, 9),
  @("Islandia", 1486, 69, 428, 1054, 11, 0, 4),
  @("Argentina", 1451, 0, 280, 1127, 86, 1, 44),
  @("Colombia", 1406, 0, 85, 1289, 50, 0, 32),
  @("Catar", 1325, 0, 109, 1213, 37, 0, 3),
  @("Argelia", 1320, 69, 90, 1078, 46, 22, 152),
  @("Singapur", 1309, 120, 320, 983, 25, 0, 6),
  @("Ucrania", 1251, 26, 25, 1194, 16, 0, 32),
  @("Croacia", 1182, 56, 125, 1042, 39, 3, 15),
  @("Estonia", 1097, 58, 62, 1020, 17, 2, 15),
  @("Egipto", 1070, 0, 241, 758, 0, 0, 71),
  @("Nueva Zelanda", 1039, 89, 156, 882, 3, 0, 1),
  @("Eslovenia", 997, 20, 79, 890, 31, 6, 28),
  @("Marruecos", 990, 71, 71, 850, 1, 10, 69),
  @("Irak", 961, 83, 279, 621, 0, 5, 61),
  @("Hong Kong", 890, 28, 206, 680, 8, 0, 4),
  @("Moldavia", 864, 112, 30, 819, 80, 3, 15),
  @("Armenia", 822, 52, 57, 758, 30, 0, 7),
  @("Lituania", 811, 40, 7, 791, 11, 2, 13),
  @("Hungria", 733, 55, 66, 633, 17, 2, 34),
  @("Crucero", 712, 0, 619, 82, 10, 0, 11),
  @("Barein", 698, 10, 427, 267, 3, 0, 4),
  @("Bosnia y Herzegovina", 654, 30, 30, 601, 4, 2, 23),
  @("Azerbaiyan", 584, 63, 32, 545, 17, 2, 7),
  @("Kazajistan", 584, 53, 42, 536, 6, 1, 6),
  @("Bielorrusia", 562, 122, 52, 502, 11, 3, 8),
  @("Kuwait", 556, 77, 99, 456, 17, 0, 1),
  @("Camerun", 555, 0, 17, 529, 0, 0, 9),
  @("Republica de Macedonia", 555, 72, 23, 514, 15, 1, 18),
  @("Tunez", 553, 0, 5, 529, 39, 1, 19),
  @("Letonia", 533, 24, 1, 531, 4, 0, 1),
  @("Bulgaria", 531, 28, 37, 474, 22, 3, 20),
  @("Libano", 527, 7, 54, 455, 28, 1, 18),
  @("Principado de Andorra", 501, 35, 26, 457, 12, 1, 18),
  @("Eslovaquia", 485, 14, 10, 474, 3, 0, 1),
  @("Republica de Chipre", 446, 20, 37, 400, 11, 0, 9),
  @("Costa Rica", 435, 0, 13, 420, 13, 0, 2),
  @("Uruguay", 400, 0, 93, 302, 13, 0, 5),
  @("Taiwan", 363, 8, 54, 304, 0, 0, 5),
  @("Albania", 361, 28, 104, 237, 7, 0, 20),
  @("Afganistan", 349, 40, 15, 327, 0, 0, 7),
  @("Burkina Faso", 345, 27, 90, 238, 0, 1, 17),
  @("Jordania", 345, 22, 110, 230, 5, 0, 5),
  @("Reunion", 344, 10, 40, 304, 4, 0, 0),
  @("Cuba", 320, 32, 15, 297, 11, 2, 8),
  @("Uzbekistan", 298, 32, 30, 266, 8, 0, 2),
  @("Oman", 298, 21, 61, 235, 3, 0, 2),
  @("Honduras", 268, 4, 6, 240, 10, 7, 22),
  @("San Marino", 266, 7, 35, 199, 14, 0, 32),
  @("Costa de Marfil", 245, 0, 25, 219, 0, 0, 1),
  @("Vietnam", 241, 1, 90, 151, 3, 0, 0),
  @("Estado de Palestina", 234, 17, 23, 210, 0, 0, 1),
  @("Malta", 227, 14, 5, 222, 3, 0, 0),
  @("Mauricio", 227, 31, 7, 213, 1, 0, 7),
  @("Nigeria", 224, 10, 27, 192, 2, 1, 5),
  @("Senegal", 222, 3, 82, 138, 1, 0, 2),
  @("Montenegro", 214, 13, 1, 211, 4, 0, 2),
  @("Ghana", 214, 9, 31, 178, 2, 0, 5),
  @("Islas Feroe", 181, 0, 99, 82, 1, 0, 0),
  @("Sri Lanka", 175, 9, 33, 137, 5, 0, 5),
  @("Georgia", 174, 12, 36, 136, 6, 1, 2),
  @("Bolivia", 157, 18, 2, 145, 3, 0, 10),
  @("Venezuela", 155, 0, 52, 96, 6, 0, 7),
  @("Consejo Danes para los Refugiados", 154, 0, 3, 133, 0, 0, 18),
  @("Martinica", 149, 4, 50, 95, 21, 1, 4),
  @("Kirguistan", 147, 3, 9, 137, 5, 0, 1),
  @("Niger", 144, 0, 0, 136, 0, 0, 8),
  @("Kenia", 142, 16, 4, 134, 2, 0, 4),
  @("Brunei", 135, 0, 73, 61, 3, 0, 1),
  @("Mayotte", 134, 0, 14, 118, 3, 0, 2),
  @("Guadalupe", 134, 0, 24, 103, 14, 0, 7),
  @("Isla de Man", 127, 1, 0, 126, 0, 0, 1),
  @("Guinea", 121, 10, 5, 116, 0, 0, 0),
  @("Camboya", 114, 0, 50, 64, 1, 0, 0),
  @("Trinidad yTobago", 104, 1, 1, 96, 0, 1, 7),
  @("Paraguay", 104, 8, 12, 89, 2, 0, 3),
  @("Gibraltar", 103, 5, 52, 51, 0, 0, 0),
  @("Ruanda", 102, 0, 4, 98, 0, 0, 0),
  @("Banglades", 88, 18, 33, 46, 1, 1, 9),
  @("Liechtenstein", 77, 0, 0, 76, 0, 0, 1),
  @("Monaco", 73, 7, 3, 69, 2, 0, 1),
  @("Madagascar", 72, 2, 2, 70, 6, 0, 0),
  @("Aruba", 64, 0, 1, 63, 0, 0, 0),
  @("El Salvador", 62, 6, 2, 57, 4, 0, 3),
  @("Guatemala", 61, 0, 15, 44, 1, 0, 2),
  @("Guayana Francesa", 61, 0, 22, 39, 1, 0, 0),
  @("Republica de Yibuti", 59, 9, 9, 50, 0, 0, 0),
  @("Barbados", 56, 4, 6, 49, 4, 1, 1),
  @("Jamaica", 55, 2, 7, 45, 0, 0, 3),
  @("Uganda", 48, 0, 0, 48, 0, 0, 0),
  @("Mali", 45, 4, 1, 39, 0, 2, 5),
  @("Congo", 45, 23, 2, 38, 0, 3, 5),
  @("Macao", 44, 0, 10, 34, 0, 0, 0),
  @("Togo", 44, 3, 20, 21, 0, 0, 3),
  @("Etiopia", 43, 5, 4, 37, 1, 2, 2),
  @("Polinesia Francesa", 40, 0, 0, 40, 1, 0, 0),
  @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
  @("Zambia", 39, 0, 3, 35, 0, 0, 1),
  @("Bermudas", 37, 0, 14, 23, 0, 0, 0),
  @("Islas Caimanes", 35, 0, 1, 33, 0, 0, 1),
  @("Guam", 32, 0, 0, 31, 0, 0, 1),
  @("Eritrea", 29, 0, 0, 29, 0, 0, 0),
  @("San Martin (Parte Francesa)", 29, 5, 7, 20, 6, 0, 2),
  @("Bahamas", 28, 0, 0, 24, 1, 0, 4),
  @("San Martin (Parte Holandesa)", 25, 2, 6, 15, 0, 0, 4),
  @("Guyana", 24, 1, 0, 20, 0, 0, 4),
  @("Benin", 22, 6, 5, 17, 0, 0, 0),
  @("Birmania", 21, 0, 0, 20, 0, 0, 1),
  @("Gabon", 21, 0, 1, 19, 0, 0, 1),
  @("Haiti", 21, 1, 1, 19, 0, 1, 1),
  @("Tanzania", 20, 0, 3, 16, 0, 0, 1),
  @("Siria", 19, 3, 2, 15, 0, 0, 2),
  @("Maldivas", 19, 0, 13, 6, 0, 0, 0),
  @("Guinea-Bisau", 18, 0, 0, 18, 0, 0, 0),
  @("Libia", 18, 0, 0, 17, 0, 0, 1),
  @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
  @("Nueva Caledonia", 17, 0, 1, 16, 0, 0, 0),
  @("Guinea Ecuatorial", 16, 0, 1, 15, 0, 0, 0),
  @("Namibia", 16, 2, 3, 13, 0, 0, 0),
  @("Antigua y Barbuda", 15, 0, 0, 15, 1, 0, 0),
  @("Dominica", 14, 0, 0, 14, 0, 0, 0),
  @("Santa Lucia", 14, 0, 1, 13, 0, 0, 0),
  @("Mongolia", 14, 0, 2, 12, 0, 0, 0),
  @("Liberia", 13, 3, 3, 7, 0, 2, 3),
  @("Fiyi", 12, 0, 0, 12, 0, 0, 0),
  @("Granada", 12, 0, 0, 12, 2, 0, 0),
  @("Sudan", 12, 2, 2, 8, 0, 0, 2),
  @("Laos", 11, 1, 0, 11, 0, 0, 0),
  @("Groenlandia", 11, 0, 3, 8, 0, 0, 0),
  @("Curazao", 11, 0, 5, 5, 0, 0, 1),
  @("San Cristobal y Nieves", 10, 1, 0, 10, 0, 0, 0),
  @("Seychelles", 10, 0, 0, 10, 0, 0, 0),
  @("Mozambique", 10, 0, 1, 9, 0, 0, 0),
  @("Surinam", 10, 0, 0, 9, 0, 0, 1),
  @("Angola", 10, 0, 2, 6, 0, 0, 2),
  @("Republica del Chad", 9, 0, 0, 9, 0, 0, 0),
  @("Suazilandia", 9, 0, 0, 9, 0, 0, 0),
  @("Nepal", 9, 0, 1, 8, 0, 0, 0),
  @("Zimbabue", 9, 0, 0, 8, 0, 0, 1),
  @("Montserrat", 9, 0, 0, 7, 0, 0, 2),
  @("Republica de Africa Central", 8, 0, 0, 8, 0, 0, 0),
  @("Santa Sede", 7, 0, 0, 7, 0, 0, 0),
  @("San Vicente y las Granadinas", 7, 0, 1, 6, 0, 0, 0),
  @("Cabo Verde", 7, 0, 0, 6, 0, 0, 1),
  @("Somalia", 7, 0, 1, 6, 0, 0, 0),
  @("Sierra Leona", 6, 2, 0, 6, 0, 0, 0),
  @("Botsuana", 6, 2, 0, 5, 0, 0, 1),
  @("San Bartolome", 6, 0, 1, 5, 0, 0, 0),
  @("Mauritania", 6, 0, 2, 3, 0, 0, 1),
  @("Belice", 5, 1, 0, 5, 1, 0, 0),
  @("Nicaragua", 5, 0, 0, 4, 0, 0, 1),
  @("Islas Turcas y Caicos", 5, 0, 0, 4, 0, 1, 1),
  @("Butan", 5, 0, 2, 3, 0, 0, 0),
  @("Sahara Occidental", 4, 0, 0, 4, 0, 0, 0),
  @("Malaui", 4, 0, 0, 4, 0, 0, 0),
  @("Gambia", 4, 0, 2, 1, 0, 0, 1),
  @("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0),
  @("Anguila", 3, 0, 0, 3, 0, 0, 0),
  @("Burundi", 3, 0, 0, 3, 0, 0, 0),
  @("Bonaire, San Eustaquio y Saba", 2, 0, 0, 2, 0, 0, 0),
  @("Sudan del Sur", 1, 1, 0, 1, 0, 0, 0),
  @("Islas Malvinas", 1, 0, 0, 1, 0, 0, 0),
  @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
  @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
  @("San Pedro y Miquelon", 1, 1, 0, 1, 0, 0, 0)
)

$r = 4
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
  $r = $r + 1
}
